$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.234.31'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.54%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.830.32'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.68%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '237.50'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.16%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6091'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -3.40%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.19%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07115'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -4.85%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2829'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.72%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.01'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -4.25%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07649'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.18%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.841.64'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.15%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.824'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -3.28%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000009997'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -2.09%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.069.75'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.89%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '79.77'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.87%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.995'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -4.58%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '29.209.71'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.57%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '230.15'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.30%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.84'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -4.10%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.16%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.056'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -4.87%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.004'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.40%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '155.44'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -2.03%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.110'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -4.63%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1299'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -3.97%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -3.88%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.06826'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +3.66%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.478'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +2.94%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.458'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -2.03%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.848'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -5.94%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.129'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.08%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.739'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -5.54%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6600'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -5.48%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.555'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.85%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.232.82'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -1.39%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.762'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -1.92%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -4.76%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.604'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -2.88%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9366'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.42%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.986.03'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.48%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.02%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '63.62'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -2.94%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.00000000117'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.44%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.635'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -4.87%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.573'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -5.19%  '
$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.1090'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -5.02%  '
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'Aptos'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.553'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -7.28%  '
